$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F ("dSF") values are being repulled/recalculated for several rows.
# Apply the updated literal values as captured by the diff.
$updates = @{
    4  = -4
    5  = -4
    10 = -5
    11 = -3
    13 = 1
    15 = 1
    17 = -4
    19 = -2
    20 = -1
    21 = 15
    22 = 0
    23 = -1
    25 = -2
    26 = 14
    27 = 3
    28 = -5
    32 = -6
    35 = -3
    36 = 1
    38 = -2
    39 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
